# ----------------------------------------------------------------------------
# Applies the "Add files via upload" edit:
#  1. Update every "Date Placeholder" field (slide master + each custom layout)
#     from 6/4/2020 -> 6/18/2020.
#  2. On the "BLOCK DAIGRAM" slide: fix the title text/spacing & grow the
#     textbox, then add a new "LDR" textbox label.
#  3. Append a new, blank slide 13 at the end of the deck.
# ----------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# 1) Date placeholders: Slide Master + every Custom Layout that has one.
# -----------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.HasTextFrame) {
      $tf = $sh.TextFrame
      if ($tf.HasText) {
        $txt = $tf.TextRange.Text
        if ($txt -eq "6/4/2020") {
          $tf.TextRange.Text = "6/18/2020"
        }
      }
    }
  }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
  $lay = $layouts.Item($li)
  Update-DatePlaceholder $lay.Shapes
}

# -----------------------------------------------------------------------
# 2) Locate the "BLOCK DAIGRAM" slide/shape and fix it up.
# -----------------------------------------------------------------------
$blockSlide = $null
$blockShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
  $s = $p.Slides.Item($i)
  for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $sh = $s.Shapes.Item($j)
    if ($sh.HasTextFrame) {
      if ($sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -match "BLOCK DAIGRAM") {
          $blockSlide = $s
          $blockShape = $sh
        }
      }
    }
  }
}

# Grow the title textbox to fit the (now two-line) heading.
$blockShape.Height = 1077218 / 12700

$titleRange = $blockShape.TextFrame.TextRange
# Replace "BLOCK DAIGRAM " (misspelled, trailing space) with the corrected,
# double-spaced "BLOCK  DIAGRAM" while leaving the 12 leading spaces intact.
$fullText = $titleRange.Text
$prefixLen = $fullText.Length - "BLOCK DAIGRAM ".Length
$oldHeading = $titleRange.Characters($prefixLen + 1, "BLOCK DAIGRAM ".Length)
$oldHeading.Text = "BLOCK  DIAGRAM"

# Add a trailing blank (centered) paragraph, matching the authored slide.
$titleRange.InsertAfter([char]13 + " ") | Out-Null

# Add the new "LDR" label textbox near the LDR sensor in the diagram.
$ldrBox = $blockSlide.Shapes.AddTextbox(1, 4519782 / 12700, 5930903 / 12700, 1235479 / 12700, 369332 / 12700)
$ldrBox.TextFrame.WordWrap = -1
$ldrBox.TextFrame.AutoSize = 1
$ldrBox.Fill.Visible = 0
$ldrBox.TextFrame.TextRange.Text = "LDR"

# -----------------------------------------------------------------------
# 3) Append a new, blank slide at the end of the deck (slide 13).
# -----------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)
